$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F2").Value = -5
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -7
